$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Straightforward per-row Price (D) / Volume(1h) (E) updates
$ws.Range("D2").Value = "98.190.54"
$ws.Range("E2").Value = "  +3.27%  "
$ws.Range("D3").Value = "3.624.11"
$ws.Range("E3").Value = "  +2.13%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.07"
$ws.Range("E5").Value = "  +4.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "659.87"
$ws.Range("E6").Value = "  +2.05%  "
$ws.Range("E7").Value = "  +20.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.419"
$ws.Range("E8").Value = "  +6.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.08"
$ws.Range("E9").Value = "  +10.18%  "
$ws.Range("E10").Value = "  -0.14%  "
$ws.Range("D11").Value = "3.620.72"
$ws.Range("E11").Value = "  +2.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.12"
$ws.Range("E12").Value = "  +5.61%  "
$ws.Range("E13").Value = "  +2.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.50"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "4.298.10"
$ws.Range("E15").Value = "  +1.60%  "
$ws.Range("D16").Value = "97.892.69"
$ws.Range("E16").Value = "  +3.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000261"
$ws.Range("E17").Value = "  +4.29%  "
$ws.Range("D18").Value = "3.618.40"
$ws.Range("E18").Value = "  +1.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.10"
$ws.Range("E19").Value = "  +3.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.92"
$ws.Range("E20").Value = "  +2.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.12"
$ws.Range("E21").Value = "  +3.16%  "
$ws.Range("E22").Value = "  +14.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.53"
$ws.Range("E23").Value = "  +2.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "515.17"
$ws.Range("E24").Value = "  +2.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000208"
$ws.Range("E25").Value = "  +9.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.93"
$ws.Range("E26").Value = "  +6.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "99.84"
$ws.Range("E27").Value = "  +6.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.02"
$ws.Range("E28").Value = "  +5.85%  "
$ws.Range("D29").Value = "3.818.81"
$ws.Range("E29").Value = "  +2.09%  "
$ws.Range("E30").Value = "  +12.61%  "
$ws.Range("E31").Value = "  +2.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.80"
$ws.Range("E32").Value = "  +6.09%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  +5.63%  "
$ws.Range("E35").Value = "  -0.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.92"
$ws.Range("E36").Value = "  +1.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.89"
$ws.Range("E37").Value = "  +9.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.573"
$ws.Range("E38").Value = "  +4.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "612.79"
$ws.Range("E39").Value = "  +10.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.64"
$ws.Range("E40").Value = "  +10.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.00"
$ws.Range("E41").Value = "  +15.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.154"
$ws.Range("E42").Value = "  +3.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.00"
$ws.Range("E45").Value = "  +7.95%  "
$ws.Range("E46").Value = "  +8.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.30"
$ws.Range("E47").Value = "  +2.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.66"
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.63"
$ws.Range("E49").Value = "  +8.54%  "

# Row 43 <-> 44 content swap (with updated values)
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.931"
$ws.Range("E43").Value = "  +4.32%  "

$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  -0.02%  "

# Row 50 <-> 51 content swap (with updated values)
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.02"
$ws.Range("E50").Value = "  -3.21%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.399"
$ws.Range("E51").Value = "  +37.17%  "
